$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 12.307962129356422
$ws.Range("C2").Value = 11.27668584925304
$ws.Range("D2").Value = 12.182262788502646
$ws.Range("E2").Value = 12.296540643036979

$ws.Range("B3").Value = 12.405401488311036
$ws.Range("C3").Value = 10.67466473685695
$ws.Range("D3").Value = 13.123817133311299
$ws.Range("E3").Value = 10.630007628015582

$ws.Range("B1:E3").Select()
